$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value. All values are written as literal
# text (matching the source file's inlineStr cells) by forcing a text number
# format before assignment, then clearing formatting so the cell style/format
# stays at the workbook default (as in the original file).
$updates = [ordered]@{
    'D2' = '37.510.80'
    'E2' = '  +2.57%  '
    'D3' = '2.079.54'
    'E3' = '  +3.74%  '
    'E4' = '  -0.03%  '
    'D5' = '235.14'
    'E5' = '  -0.72%  '
    'E6' = '  +3.23%  '
    'D7' = '58.19'
    'E7' = '  +5.99%  '
    'E8' = '  +0.03%  '
    'E9' = '  +3.71%  '
    'D10' = '59.10'
    'E10' = '  +1.31%  '
    'D11' = '0.0763'
    'E11' = '  +2.19%  '
    'E12' = '  +3.93%  '
    'D13' = '2.386.45'
    'E13' = '  +3.71%  '
    'D14' = '14.59'
    'E14' = '  +2.80%  '
    'D15' = '21.12'
    'E15' = '  +4.51%  '
    'E16' = '  +3.16%  '
    'E17' = '  +2.33%  '
    'D18' = '2.079.02'
    'E18' = '  +3.61%  '
    'D19' = '37.477.63'
    'E19' = '  +2.64%  '
    'D20' = '6.27'
    'E20' = '  +18.07%  '
    'D21' = '70.12'
    'E21' = '  +3.41%  '
    'E22' = '  +1.56%  '
    'D23' = '226.74'
    'E23' = '  +2.27%  '
    'D24' = '0.999'
    'E24' = '  -0.20%  '
    'E25' = '  +2.63%  '
    'E26' = '  +0.58%  '
    'D27' = '166.88'
    'E27' = '  +2.63%  '
    'D28' = '1.52'
    'E28' = '  +11.92%  '
    'D29' = '9.04'
    'E29' = '  +4.35%  '
    'D30' = '19.30'
    'E30' = '  +2.80%  '
    'E31' = '  -0.73%  '
    'D32' = '0.118'
    'E32' = '  +1.22%  '
    'D33' = '4.55'
    'E33' = '  +3.79%  '
    'E34' = '  +3.37%  '
    'D35' = '2.59'
    'E35' = '  +6.24%  '
    'E36' = '  +7.49%  '
    'E37' = '  -0.06%  '
    'E38' = '  +0.30%  '
    'B39' = 'THORChain'
    'C39' = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
    'D39' = '5.89'
    'E39' = '  +3.35%  '
    'B40' = 'WEMIXToken'
    'C40' = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
    'D40' = '1.77'
    'E40' = '  +0.51%  '
    'D41' = '4.65'
    'E41' = '  +21.47%  '
    'E42' = '  -1.09%  '
    'D43' = '0.0958'
    'E43' = '  +3.54%  '
    'D44' = '1.476.55'
    'E44' = '  +1.66%  '
    'B45' = 'Aave'
    'C45' = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    'D45' = '96.15'
    'E45' = '  +6.19%  '
    'B46' = 'TrustWalletToken'
    'C46' = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    'D46' = '1.18'
    'E46' = '  +7.08%  '
    'E47' = '  +4.89%  '
    'D48' = '15.87'
    'E48' = '  +3.95%  '
    'E49' = '  +3.93%  '
    'D50' = '7.29'
    'E50' = '  +6.09%  '
    'E51' = '  +1.79%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = '@'
    $cell.Value = $updates[$ref]
    $cell.ClearFormats()
}

